# Auto-generated edit script: refresh market-board derived values
# across the per-job Leve profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# ALC!row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 35910.066
$ws.Range("I64").Value = 168866.67
$ws.Range("J64").Value = 2670.9167
$ws.Range("K64").Value = 168866.67
$ws.Range("L64").Value = 2670.9167
$ws.Range("M64").Value = -168618.67
$ws.Range("N64").Value = -3166.9167

# ALC!row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 35910.066
$ws.Range("I67").Value = 168866.67
$ws.Range("J67").Value = 2670.9167
$ws.Range("K67").Value = 168866.67
$ws.Range("L67").Value = 2670.9167
$ws.Range("M67").Value = -168008.67
$ws.Range("N67").Value = -4386.9167

# ALC!row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1455.8889
$ws.Range("I70").Value = 1333.3334
$ws.Range("J70").Value = 1517.1666
$ws.Range("K70").Value = 4000.0002
$ws.Range("L70").Value = 4551.4998
$ws.Range("M70").Value = -3730.0002
$ws.Range("N70").Value = -5091.4998

# ALC!row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1455.8889
$ws.Range("I73").Value = 1333.3334
$ws.Range("J73").Value = 1517.1666
$ws.Range("K73").Value = 4000.0002
$ws.Range("L73").Value = 4551.4998
$ws.Range("M73").Value = -3064.0002
$ws.Range("N73").Value = -6423.4998

# ALC!row 93
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 46186.25
$ws.Range("J93").Value = 46186.25
$ws.Range("L93").Value = 46186.25
$ws.Range("N93").Value = -51178.25

# ALC!row 95
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 39742
$ws.Range("J95").Value = 39742
$ws.Range("L95").Value = 39742
$ws.Range("N95").Value = -45234

# ALC!row 97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 74493.336
$ws.Range("J97").Value = 74493.336
$ws.Range("L97").Value = 223480.008
$ws.Range("N97").Value = -224472.008

# ALC!row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1512.6
$ws.Range("I100").Value = 1298.7778
$ws.Range("J100").Value = 1833.3334
$ws.Range("K100").Value = 1298.7778
$ws.Range("L100").Value = 1833.3334
$ws.Range("M100").Value = -757.7778000000001
$ws.Range("N100").Value = -2915.3334

# ALC!row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 455.7857
$ws.Range("J101").Value = 1350
$ws.Range("L101").Value = 4050
$ws.Range("N101").Value = -7294

# ALC!row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2616.9714
$ws.Range("I113").Value = 2525.5
$ws.Range("J113").Value = 2644.074
$ws.Range("K113").Value = 2525.5
$ws.Range("L113").Value = 2644.074
$ws.Range("M113").Value = 728.5
$ws.Range("N113").Value = -9152.074000000001

# ALC!row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6120.5
$ws.Range("I116").Value = 4148.1816
$ws.Range("J116").Value = 7566.8667
$ws.Range("K116").Value = 4148.1816
$ws.Range("L116").Value = 7566.8667
$ws.Range("M116").Value = -706.1815999999999
$ws.Range("N116").Value = -14450.8667

# ALC!row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 22642.783
$ws.Range("I132").Value = 3495.8108
$ws.Range("K132").Value = 10487.4324
$ws.Range("M132").Value = -7957.432400000002

# ARM!row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1331.8125
$ws.Range("I97").Value = 1129.8
$ws.Range("K97").Value = 1129.8
$ws.Range("M97").Value = -633.8

# BSM!row 100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 46409
$ws.Range("J100").Value = 46409
$ws.Range("L100").Value = 46409
$ws.Range("N100").Value = -48573

# CRP!row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3019.4375
$ws.Range("I62").Value = 2923
$ws.Range("J62").Value = 3437.3333
$ws.Range("K62").Value = 2923
$ws.Range("L62").Value = 3437.3333
$ws.Range("M62").Value = -2299
$ws.Range("N62").Value = -4685.3333

# CRP!row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3019.4375
$ws.Range("I65").Value = 2923
$ws.Range("J65").Value = 3437.3333
$ws.Range("K65").Value = 14615
$ws.Range("L65").Value = 17186.6665
$ws.Range("M65").Value = -11495
$ws.Range("N65").Value = -23426.6665

# CUL!row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 17850.334
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 17850.334
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 53551.00199999999
$ws.Range("N22").Value = -53889.00199999999
$ws.Range("M22").ClearContents()

# CUL!row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 17850.334
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 17850.334
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 53551.00199999999
$ws.Range("N27").Value = -53755.00199999999
$ws.Range("M27").ClearContents()

# CUL!row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1702246.9
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1702246.9
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 5106740.699999999
$ws.Range("N58").Value = -5106996.699999999
$ws.Range("M58").ClearContents()

# CUL!row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2658.8333
$ws.Range("I132").Value = 1409.7142
$ws.Range("J132").Value = 3173.1765
$ws.Range("K132").Value = 12687.4278
$ws.Range("L132").Value = 28558.5885
$ws.Range("M132").Value = -10157.4278
$ws.Range("N132").Value = -33618.5885

# GSM!row 75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 38526.2
$ws.Range("J75").Value = 38526.2
$ws.Range("L75").Value = 38526.2
$ws.Range("N75").Value = -40274.2

# GSM!row 78
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H78").Value = 38526.2
$ws.Range("J78").Value = 38526.2
$ws.Range("L78").Value = 115578.6
$ws.Range("N78").Value = -124314.6

# GSM!row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 4524.4116
$ws.Range("I97").Value = 4461.875
$ws.Range("J97").Value = 4580
$ws.Range("K97").Value = 4461.875
$ws.Range("L97").Value = 4580
$ws.Range("M97").Value = -3965.875
$ws.Range("N97").Value = -5572

# LTW!row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1633.8334
$ws.Range("I22").Value = 801
$ws.Range("J22").Value = 1800.4
$ws.Range("K22").Value = 801
$ws.Range("L22").Value = 1800.4
$ws.Range("M22").Value = -506
$ws.Range("N22").Value = -2390.4

# LTW!row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1633.8334
$ws.Range("I27").Value = 801
$ws.Range("J27").Value = 1800.4
$ws.Range("K27").Value = 801
$ws.Range("L27").Value = 1800.4
$ws.Range("M27").Value = -694
$ws.Range("N27").Value = -2014.4

# LTW!row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2899.9167
$ws.Range("I93").Value = 3000
$ws.Range("K93").Value = 3000
$ws.Range("M93").Value = -1752

# LTW!row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2420.9285
$ws.Range("I100").Value = 2049.125
$ws.Range("J100").Value = 2916.6667
$ws.Range("K100").Value = 2049.125
$ws.Range("L100").Value = 2916.6667
$ws.Range("M100").Value = -1508.125
$ws.Range("N100").Value = -3998.6667

# LTW!row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2151.05
$ws.Range("I122").Value = 2028.0667
$ws.Range("J122").Value = 2520
$ws.Range("K122").Value = 6084.2001
$ws.Range("L122").Value = 7560
$ws.Range("M122").Value = -3634.2001
$ws.Range("N122").Value = -12460

# WVR!row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 446
$ws.Range("I100").Value = 369.8
$ws.Range("K100").Value = 739.6
$ws.Range("M100").Value = -198.6

# WVR!row 137
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 38777.69
$ws.Range("J137").Value = 38777.69
$ws.Range("L137").Value = 38777.69
$ws.Range("N137").Value = -48977.69
